# Apply the updates described by the diff: numeric "want to go" counts
# (column F) and "min ticket price" values (column G) were refreshed for
# a handful of rows across the 展览, 演出 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$wsExhibit.Range("G3").Value  = 29
$wsExhibit.Range("F7").Value  = 7630
$wsExhibit.Range("F11").Value = 8257
$wsExhibit.Range("F15").Value = 5642
$wsExhibit.Range("F17").Value = 2599
$wsExhibit.Range("F19").Value = 4593
$wsExhibit.Range("F23").Value = 31
$wsExhibit.Range("G23").Value = 39
$wsExhibit.Range("F24").Value = 522
$wsExhibit.Range("F25").Value = 3455
$wsExhibit.Range("F26").Value = 39
$wsExhibit.Range("F29").Value = 2939
$wsExhibit.Range("F31").Value = 59
$wsExhibit.Range("F34").Value = 126
$wsExhibit.Range("F35").Value = 299
$wsExhibit.Range("F36").Value = 171
$wsExhibit.Range("F37").Value = 652
$wsExhibit.Range("F44").Value = 2710

# --- 演出 (sheet2) ---
$wsShow.Range("F3").Value = 118
$wsShow.Range("F5").Value = 47
$wsShow.Range("F7").Value = 37

# --- 全部类型 (sheet4) ---
$wsAll.Range("G5").Value  = 29
$wsAll.Range("F7").Value  = 7630
$wsAll.Range("F11").Value = 8257
$wsAll.Range("F15").Value = 5642
$wsAll.Range("F17").Value = 2599
$wsAll.Range("F19").Value = 4593
$wsAll.Range("F23").Value = 31
$wsAll.Range("G23").Value = 39
$wsAll.Range("F24").Value = 118
$wsAll.Range("F25").Value = 522
$wsAll.Range("F26").Value = 3455
$wsAll.Range("F29").Value = 2939
$wsAll.Range("F31").Value = 126
$wsAll.Range("F32").Value = 299
$wsAll.Range("F33").Value = 47
$wsAll.Range("F34").Value = 652
$wsAll.Range("F38").Value = 37
$wsAll.Range("F43").Value = 2710
